# Scheduled market-data refresh: updates computed price/profit columns
# (H..N) for the affected leve rows across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 5660.6665
$ws.Range("J40").Value = 4742.5
$ws.Range("L40").Value = 4742.5
$ws.Range("N40").Value = -5092.5
# Row 88: The Grave of Hemlock Groves | Growth Formula Zeta
$ws.Range("H88").Value = 2555
$ws.Range("I88").Value = 3199.6667
$ws.Range("J88").Value = 2416.8572
$ws.Range("K88").Value = 3199.6667
$ws.Range("L88").Value = 2416.8572
$ws.Range("M88").Value = -2793.6667
$ws.Range("N88").Value = -3228.8572
# Row 91: Dappling the Highlands (L) | Growth Formula Zeta
$ws.Range("H91").Value = 2555
$ws.Range("I91").Value = 3199.6667
$ws.Range("J91").Value = 2416.8572
$ws.Range("K91").Value = 3199.6667
$ws.Range("L91").Value = 2416.8572
$ws.Range("M91").Value = -1795.6667
$ws.Range("N91").Value = -5224.8572
# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 2805143
$ws.Range("I98").Value = 3368871.8
$ws.Range("K98").Value = 3368871.8
$ws.Range("M98").Value = -3367373.8
# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 2805143
$ws.Range("I122").Value = 3368871.8
$ws.Range("K122").Value = 10106615.4
$ws.Range("M122").Value = -10104165.4
# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 3553.7354
$ws.Range("I132").Value = 2937.6428
$ws.Range("K132").Value = 8812.928400000001
$ws.Range("M132").Value = -6282.928400000001
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3658.9744
$ws.Range("J138").Value = 3998.375
$ws.Range("L138").Value = 11995.125
$ws.Range("N138").Value = -22275.125

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 1848
$ws.Range("I2").Value = 2003.6666
$ws.Range("K2").Value = 2003.6666
$ws.Range("M2").Value = -1890.6666
# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 1848
$ws.Range("I116").Value = 2003.6666
$ws.Range("K116").Value = 2003.6666
$ws.Range("M116").Value = 290.3334
# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 3284.2942
$ws.Range("I122").Value = 3402.0625
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 10206.1875
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -7756.1875
$ws.Range("N122").Value = -9100
# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 2327.074
$ws.Range("I132").Value = 2288.8333
$ws.Range("J132").Value = 2633
$ws.Range("K132").Value = 6866.499899999999
$ws.Range("L132").Value = 7899
$ws.Range("M132").Value = -4336.499899999999
$ws.Range("N132").Value = -12959

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 1848
$ws.Range("I3").Value = 2003.6666
$ws.Range("K3").Value = 2003.6666
$ws.Range("M3").Value = -1889.6666
# Row 138: Bladewinner | Titanium Gold Greatsword
$ws.Range("H138").Value = 63150.59
$ws.Range("J138").Value = 63150.59
$ws.Range("L138").Value = 63150.59
$ws.Range("N138").Value = -73430.59

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 251839.08
$ws.Range("I31").Value = 358383.97
$ws.Range("J31").Value = 3234.3333
$ws.Range("K31").Value = 358383.97
$ws.Range("L31").Value = 3234.3333
$ws.Range("M31").Value = -358088.97
$ws.Range("N31").Value = -3824.3333
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 251839.08
$ws.Range("I34").Value = 358383.97
$ws.Range("J34").Value = 3234.3333
$ws.Range("K34").Value = 358383.97
$ws.Range("L34").Value = 3234.3333
$ws.Range("M34").Value = -358181.97
$ws.Range("N34").Value = -3638.3333
# Row 68: Do You Even String Bow | Holy Cedar Composite Bow
$ws.Range("H68").Value = 41779
$ws.Range("J68").Value = 41779
$ws.Range("L68").Value = 41779
$ws.Range("N68").Value = -43277
# Row 71: Win One Bow, Get Three Free (L) | Holy Cedar Composite Bow
$ws.Range("H71").Value = 41779
$ws.Range("J71").Value = 41779
$ws.Range("L71").Value = 125337
$ws.Range("N71").Value = -132825
# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 428440.78
$ws.Range("I99").Value = 840212.0600000001
$ws.Range("K99").Value = 840212.0600000001
$ws.Range("M99").Value = -838714.0600000001
# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 4749.3613
$ws.Range("I107").Value = 734.5454999999999
$ws.Range("K107").Value = 734.5454999999999
$ws.Range("M107").Value = 1185.4545
# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 428440.78
$ws.Range("I126").Value = 840212.0600000001
$ws.Range("K126").Value = 2520636.18
$ws.Range("M126").Value = -2518166.18
# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 2412.138
$ws.Range("I132").Value = 2386.24
$ws.Range("K132").Value = 7158.719999999999
$ws.Range("M132").Value = -4628.719999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 987.6
$ws.Range("I5").Value = 984.5
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 2953.5
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -2841.5
$ws.Range("N5").Value = -3224
# Row 9: Jack of All Plates | Jack-o'-lantern
$ws.Range("H9").Value = 814.7
$ws.Range("I9").Value = 257.25
$ws.Range("K9").Value = 771.75
$ws.Range("M9").Value = -547.75
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 9092344
$ws.Range("I131").Value = 71429416
$ws.Range("K131").Value = 214288248
$ws.Range("M131").Value = -214283208
# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 987.6
$ws.Range("I135").Value = 984.5
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 8860.5
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6325.5
$ws.Range("N135").Value = -14070
# Row 137: Creative Chocolate | Gateau au Chocolat
$ws.Range("H137").Value = 1625
$ws.Range("J137").Value = 2000
$ws.Range("L137").Value = 6000
$ws.Range("N137").Value = -16200
# Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws.Range("H140").Value = 16644.934
$ws.Range("I140").Value = 41542.8
$ws.Range("K140").Value = 124628.4
$ws.Range("M140").Value = -119448.4

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell | Hardsilver Ingot
$ws.Range("H80").Value = 3692.9697
$ws.Range("I80").Value = 2219.5
$ws.Range("K80").Value = 2219.5
$ws.Range("M80").Value = -1221.5
# Row 83: With a Noise That Reaches Heaven (L) | Hardsilver Ingot
$ws.Range("H83").Value = 3692.9697
$ws.Range("I83").Value = 2219.5
$ws.Range("K83").Value = 11097.5
$ws.Range("M83").Value = -6105.5
# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 2936.9412
$ws.Range("I97").Value = 2462.3
$ws.Range("K97").Value = 2462.3
$ws.Range("M97").Value = -1966.3
# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 2105.4707
$ws.Range("I107").Value = 2158.6667
$ws.Range("J107").Value = 1977.8
$ws.Range("K107").Value = 2158.6667
$ws.Range("L107").Value = 1977.8
$ws.Range("M107").Value = -238.6667000000002
$ws.Range("N107").Value = -5817.8
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 6875.7207
$ws.Range("I122").Value = 6762.8887
$ws.Range("K122").Value = 20288.6661
$ws.Range("M122").Value = -17838.6661
# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 10294.625
$ws.Range("I126").Value = 13939.4
$ws.Range("K126").Value = 41818.2
$ws.Range("M126").Value = -39348.2
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 23356.588
$ws.Range("I132").Value = 26482.068
$ws.Range("K132").Value = 79446.204
$ws.Range("M132").Value = -76916.204

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 7045.294
$ws.Range("I7").Value = 6978.3335
$ws.Range("K7").Value = 6978.3335
$ws.Range("M7").Value = -6866.3335
# Row 68: You Could Say It's a Moving Target | Wyvern Leather
$ws.Range("H68").Value = 8012.174
$ws.Range("I68").Value = 8914.549999999999
$ws.Range("J68").Value = 1996.3334
$ws.Range("K68").Value = 8914.549999999999
$ws.Range("L68").Value = 1996.3334
$ws.Range("M68").Value = -8165.549999999999
$ws.Range("N68").Value = -3494.3334
# Row 71: They Call It Bloody Mary (L) | Wyvern Leather
$ws.Range("H71").Value = 8012.174
$ws.Range("I71").Value = 8914.549999999999
$ws.Range("J71").Value = 1996.3334
$ws.Range("K71").Value = 44572.75
$ws.Range("L71").Value = 9981.666999999999
$ws.Range("M71").Value = -40828.75
$ws.Range("N71").Value = -17469.667
# Row 93: Hide to Go Seek | Gagana Leather
$ws.Range("H93").Value = 71431700
$ws.Range("I93").Value = 2748.7778
$ws.Range("J93").Value = 200003790
$ws.Range("K93").Value = 2748.7778
$ws.Range("L93").Value = 200003790
$ws.Range("M93").Value = -1500.7778
$ws.Range("N93").Value = -200006286
# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 7045.294
$ws.Range("I126").Value = 6978.3335
$ws.Range("K126").Value = 20935.0005
$ws.Range("M126").Value = -18465.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables | Hempen Underpants
$ws.Range("H2").Value = 99999
$ws.Range("I2").Value = 99999
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 99999
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -99887
$ws.Range("N2").ClearContents()
# Row 69: Fashion Patrol | Holy Rainbow Sarouel of Casting
$ws.Range("H69").Value = 6630.25
$ws.Range("J69").Value = 6630.25
$ws.Range("L69").Value = 6630.25
$ws.Range("N69").Value = -8128.25
# Row 72: Dress Code Violation (L) | Holy Rainbow Sarouel of Casting
$ws.Range("H72").Value = 6630.25
$ws.Range("J72").Value = 6630.25
$ws.Range("L72").Value = 19890.75
$ws.Range("N72").Value = -27378.75
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 4510.125
$ws.Range("I81").Value = 5484.8667
$ws.Range("K81").Value = 10969.7334
$ws.Range("M81").Value = -9908.733399999999
# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 4510.125
$ws.Range("I84").Value = 5484.8667
$ws.Range("K84").Value = 54848.66699999999
$ws.Range("M84").Value = -49544.66699999999
# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 40408.184
$ws.Range("I96").Value = 65527.5
$ws.Range("J96").Value = 3871
$ws.Range("K96").Value = 65527.5
$ws.Range("L96").Value = 3871
$ws.Range("M96").Value = -64154.5
$ws.Range("N96").Value = -6617
# Row 107: Flax Wax | Bright Linen Yarn
$ws.Range("H107").Value = 597.58826
$ws.Range("J107").Value = 733
$ws.Range("L107").Value = 2199
$ws.Range("N107").Value = -6039
# Row 112: Hair Do No Harm | Iridescent Hat of Healing
$ws.Range("H112").Value = 53615.445
$ws.Range("J112").Value = 53615.445
$ws.Range("L112").Value = 53615.445
$ws.Range("N112").Value = -56569.445
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1229.9231
$ws.Range("I132").Value = 1219.1
$ws.Range("K132").Value = 3657.3
$ws.Range("M132").Value = -1127.3
